$d = $word.ActiveDocument

$replacements = @(
    @("852÷2=", "752÷3="),
    @("525÷5=", "346÷4="),
    @("890÷6=", "722÷8="),
    @("486÷5=", "587÷8="),
    @("647÷8=", "832÷9="),
    @("898÷9=", "159÷2="),
    @("245÷5=", "249÷9="),
    @("330÷4=", "865÷9="),
    @("706÷5=", "629÷8="),
    @("586÷5=", "880÷4="),
    @("265÷9=", "418÷3="),
    @("439÷4=", "978÷5="),
    @("346÷8=", "375÷7="),
    @("605÷6=", "633÷7="),
    @("424÷2=", "706÷4="),
    @("261÷7=", "710÷8="),
    @("829÷5=", "494÷5="),
    @("200÷6=", "509÷5="),
    @("948÷7=", "385÷5="),
    @("809÷2=", "520÷7="),
    @("206÷4=", "628÷4="),
    @("812÷8=", "494÷6="),
    @("800÷4=", "844÷9="),
    @("433÷6=", "429÷5="),
    @("826÷5=", "250÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
